$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.129888534545898
$ws.Range("B1").Value = 2.127791166305542
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.252900123596191
$ws.Range("E1").Value = 1.087620139122009
